# This workbook's scraper re-runs periodically, re-stamping the "取得日時"
# (retrieved-at) column for every listing currently on the ランサーズ sheet
# with the timestamp of the latest run: 2025-11-15 12:42:13 (JST), which
# was previously 2025-11-15 12:32:28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-15 12:42:13"

# Find the last used row (header is row 1, data starts at row 2).
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
